$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.731.05"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "3.450.24"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.19"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.03"
$ws.Range("E6").Value = "  +7.30%  "
$ws.Range("D7").Value = "3.450.55"
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("E11").Value = "  +3.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.392"
$ws.Range("E12").Value = "  +2.78%  "
$ws.Range("D13").Value = "4.038.64"
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.95"
$ws.Range("E14").Value = "  +9.79%  "
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "3.447.87"
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").Value = "61.801.72"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.25"
$ws.Range("E19").Value = "  +8.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.42"
$ws.Range("E20").Value = "  +4.27%  "
$ws.Range("E21").Value = "  +2.25%  "
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.46"
$ws.Range("E24").Value = "  +3.56%  "
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("D28").Value = "3.600.40"
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.71"
$ws.Range("E30").Value = "  +4.45%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.49"
$ws.Range("E32").Value = "  -9.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.22"
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "24.16"
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("D37").Value = "3.479.24"
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.00"
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "166.53"
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0786"
$ws.Range("E42").Value = "  +4.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "27.07"
$ws.Range("E43").Value = "  +7.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.807"
$ws.Range("E44").Value = "  +4.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.55"
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("E46").Value = "  +3.76%  "
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("D50").Value = "2.572.26"
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.93"
$ws.Range("E51").Value = "  +2.16%  "
